# Methods-section edit: the footnote/legend paragraph under the final
# regression table currently ends with the abbreviation note for "CORs"
# ("... CORs crude odds ratio"). Add a comma after "crude odds ratio"
# and append a new abbreviation entry for "AIC" (Akaike Information
# Criterion), matching the style of the existing entries (italic
# abbreviation run followed by a plain-text expansion run).

$d = $word.ActiveDocument

$find = $d.Content.Find
$found = $find.Execute(" crude odds ratio", $false, $false, $false, $false, `
                        $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the ' crude odds ratio' text to edit."
}

$rng = $find.Parent

# Append the comma + trailing space to the existing run's text.
$rng.Text = " crude odds ratio, "

# Move to just after the text we just wrote, then insert the new
# italic "AIC" abbreviation run.
$rng.Collapse(0)
$insertStart = $rng.Start
$rng.InsertAfter("AIC Akaike Information Criterion")

# Italicize only the "AIC" portion (first 3 characters of what we just
# inserted) so it becomes its own run, matching the formatting pattern
# used by the other abbreviation/definition pairs in this paragraph.
$italicRng = $d.Range($insertStart, $insertStart + 3)
$italicRng.Font.Italic = $true
